$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.913
$ws.Range("D4").Value = -7.676
$ws.Range("E4").Value = 13.091

$ws.Range("D5").Value = -8.206999999999999

$ws.Range("C7").Value = -13.088

$ws.Range("D8").Value = -7.896000000000001

$ws.Range("E9").Value = 12.946

$ws.Range("C16").Value = -12.302
$ws.Range("D16").Value = -8.574000000000002

$ws.Range("E18").Value = 13.451
